$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new record as row 8, pushing all existing data rows
# (old rows 8-86) down by one (new rows 9-87).
$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = 10
$ws.Range("B8").Value = "Vega Modelo de Temuco"
$ws.Range("C8").Value = "La Araucanía"
$ws.Range("D8").Value = 44847
$ws.Range("E8").Value = 9
$ws.Range("F8").Value = 300000001
$ws.Range("G8").Value = "Rabanito"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 190
$ws.Range("K8").Value = 6000
$ws.Range("L8").Value = 8000
$ws.Range("M8").Value = 7000
$ws.Range("N8").Value = "$/docena de paquetes"
$ws.Range("O8").Value = "Provincia de Cautín"
$ws.Range("P8").Value = 583
$ws.Range("Q8").Value = 12
$ws.Range("R8").Value = "Hortaliza"
